$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (rows 10, 13-24) to reflect reorganized content ---
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "1341653 - Maria José Ramos Sandim"
$ws.Range("C15").Value = "1341653 - Maria José Ramos Sandim"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Range("C18").Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas e práticas ministradas em laboratório."
$ws.Range("C19").Value = "Aulas expositivas e práticas ministradas em laboratório."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOB1021 -  Física IV  (Requisito)"
$ws.Range("C23").Value = "LOB1021 -  Física IV  (Requisito)"
$ws.Range("B24").Value = "LOM3205 -  Eletromagnetismo  (Requisito)"
$ws.Range("C24").Value = "LOM3205 -  Eletromagnetismo  (Requisito)"

# --- Clear cells that no longer hold content ---
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()
$ws.Range("A24").Clear()

# --- Adjust row heights (rows 13-24) to match the new layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30

# --- Remove the trailing rows (old rows 25-27 are no longer needed) ---
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(25).Delete()

